# Fruta / hortaliza, semanal
#
# A new weekly price record was added to the "Mango" sheet for the
# Terminal Hortofrutícola Agro Chillán market. The new observation is
# inserted as row 54 (pushing the previously existing rows 54-58 down
# to rows 55-59), matching the other rows' static columns and carrying
# its own date / volume / price / origin values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54; rows 54-58 shift down to 55-59.
$ws.Rows("54").Insert()

$ws.Cells.Item(54, 1).Value2 = 7
$ws.Cells.Item(54, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(54, 3).Value2 = "Ñuble"
$ws.Cells.Item(54, 4).Value2 = 44474
$ws.Cells.Item(54, 5).Value2 = 16
$ws.Cells.Item(54, 6).Value2 = "Fruta"
$ws.Cells.Item(54, 7).Value2 = 100108
$ws.Cells.Item(54, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(54, 9).Value2 = 100108002
$ws.Cells.Item(54, 10).Value2 = "Mango"
$ws.Cells.Item(54, 11).Value2 = "Sin especificar"
$ws.Cells.Item(54, 12).Value2 = "Primera"
$ws.Cells.Item(54, 13).Value2 = 60
$ws.Cells.Item(54, 14).Value2 = 8500
$ws.Cells.Item(54, 15).Value2 = 9000
$ws.Cells.Item(54, 16).Value2 = 8750
$ws.Cells.Item(54, 17).Value2 = "`$/bandeja 4 kilos"
$ws.Cells.Item(54, 18).Value2 = "Brasil"
$ws.Cells.Item(54, 19).Value2 = 2188
$ws.Cells.Item(54, 20).Value2 = 4
